$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (shifts existing PercActivations* columns from B:F to C:G)
$ws.Columns("B").Insert()

# New header for the inserted column
$ws.Range("B1").Value = "segments"

# Move the segment-name text (currently still in column A, rows 2-20) into the
# new column B, and replace column A with a numeric 0-based index.
for ($i = 2; $i -le 20; $i++) {
    $segmentName = $ws.Range("A" + $i).Value()
    $ws.Range("B" + $i).Value = $segmentName
    $ws.Range("A" + $i).Value = $i - 2
}

# The inserted column B picked up column A's header-row formatting; the data
# rows (2-20) should have no explicit style, matching the target layout.
$ws.Range("B2:B20").ClearFormats()

# Give the new B1 header cell the same style as the other header cells.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
